$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (the "Förändrad" / last-changed date) holds the same serial
# date value (45177 = 2023-09-08) for every data row from row 2 to 264.
# Bump it to 45178 (2023-09-09) for all of them.
$ws.Range("C2:C264").Value = 45178
